# Update crypto price/volume figures per the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.128.51"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "'1.835.29"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("D4").Value = "'0.9991"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'240.49"
$ws.Range("E5").Value = "  -2.16%  "

$ws.Range("D6").Value = "'0.6841"
$ws.Range("E6").Value = "  -2.03%  "

$ws.Range("D7").Value = "'0.9996"

$ws.Range("D8").Value = "'0.3011"
$ws.Range("E8").Value = "  -1.60%  "

$ws.Range("D9").Value = "'0.07447"
$ws.Range("E9").Value = "  -3.44%  "

$ws.Range("E10").Value = "  -2.21%  "

$ws.Range("D11").Value = "'0.07663"
$ws.Range("E11").Value = "  -2.01%  "

$ws.Range("D12").Value = "'1.836.00"
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").Value = "'5.049"
$ws.Range("E13").Value = "  -1.58%  "

$ws.Range("D14").Value = "'0.6813"
$ws.Range("E14").Value = "  -0.78%  "

$ws.Range("D15").Value = "'87.53"
$ws.Range("E15").Value = "  -6.23%  "

$ws.Range("D16").Value = "'6.146"
$ws.Range("E16").Value = "  -7.63%  "

$ws.Range("D17").Value = "'29.114.36"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").Value = "'0.000008171"
$ws.Range("E18").Value = "  -1.61%  "

$ws.Range("D19").Value = "'2.081.27"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").Value = "'227.65"
$ws.Range("E20").Value = "  -5.95%  "

$ws.Range("E21").Value = "  -2.02%  "

$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "'7.399"
$ws.Range("E23").Value = "  -1.53%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").Value = "'0.1455"
$ws.Range("E25").Value = "  -3.97%  "

$ws.Range("D26").Value = "'159.96"
$ws.Range("E26").Value = "  +0.23%  "

$ws.Range("D27").Value = "'8.746"
$ws.Range("E27").Value = "  -0.97%  "

$ws.Range("D28").Value = "'18.10"
$ws.Range("E28").Value = "  -0.94%  "

$ws.Range("D29").Value = "'1.511"
$ws.Range("E29").Value = "  -2.00%  "

$ws.Range("E30").Value = "  +0.45%  "

$ws.Range("D31").Value = "'4.150"
$ws.Range("E31").Value = "  -0.95%  "

$ws.Range("D32").Value = "'1.194"
$ws.Range("E32").Value = "  -0.33%  "

$ws.Range("D33").Value = "'0.05160"
$ws.Range("E33").Value = "  +0.70%  "

$ws.Range("D34").Value = "'0.7653"
$ws.Range("E34").Value = "  -3.00%  "

$ws.Range("D35").Value = "'1.841"
$ws.Range("E35").Value = "  -1.21%  "

$ws.Range("D36").Value = "'1.133"
$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("D37").Value = "'2.674"
$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("D38").Value = "'1.303.80"
$ws.Range("E38").Value = "  -1.23%  "

$ws.Range("D39").Value = "'0.01830"
$ws.Range("E39").Value = "  -2.10%  "

$ws.Range("D40").Value = "'2.719"
$ws.Range("E40").Value = "  +0.29%  "

$ws.Range("D41").Value = "'0.9334"
$ws.Range("E41").Value = "  -1.32%  "

$ws.Range("D42").Value = "'5.805"
$ws.Range("E42").Value = "  -4.28%  "

$ws.Range("D43").Value = "'104.52"
$ws.Range("E43").Value = "  -3.06%  "

$ws.Range("D44").Value = "'0.9988"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("E45").Value = "  +0.25%  "

$ws.Range("D46").Value = "'65.21"
$ws.Range("E46").Value = "  +1.34%  "

$ws.Range("D47").Value = "'1.981.66"
$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").Value = "'0.5196"
$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").Value = "'9.516"
$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("D50").Value = "'1.769"
$ws.Range("E50").Value = "  +0.26%  "

$ws.Range("D51").Value = "'0.05931"
$ws.Range("E51").Value = "  +0.77%  "
